$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

Set-TextValue 'D2' '26.902.22'
Set-TextValue 'E2' '  +2.17%  '
Set-TextValue 'D3' '1.812.26'
Set-TextValue 'E3' '  +2.95%  '
Set-TextValue 'E4' '  +0.55%  '
Set-TextValue 'D5' '314.16'
Set-TextValue 'E5' '  +3.45%  '
Set-TextValue 'E6' '  +0.57%  '
Set-TextValue 'D7' '0.4292'
Set-TextValue 'E7' '  +0.38%  '
Set-TextValue 'D8' '0.3693'
Set-TextValue 'E8' '  +2.35%  '
Set-TextValue 'D9' '0.07242'
Set-TextValue 'E9' '  +2.69%  '
Set-TextValue 'D10' '0.8662'
Set-TextValue 'E10' '  +4.50%  '
Set-TextValue 'D11' '2.045.84'
Set-TextValue 'E11' '  +17.97%  '
Set-TextValue 'D12' '21.28'
Set-TextValue 'E12' '  +5.91%  '
Set-TextValue 'D13' '5.406'
Set-TextValue 'D14' '6.626'
Set-TextValue 'E14' '  +4.38%  '
Set-TextValue 'D15' '0.06948'
Set-TextValue 'E15' '  +1.88%  '
Set-TextValue 'D16' '80.82'
Set-TextValue 'E16' '  +2.32%  '
Set-TextValue 'D17' '1.007'
Set-TextValue 'E17' '  +0.21%  '
Set-TextValue 'D18' '0.000008925'
Set-TextValue 'E18' '  +3.36%  '
Set-TextValue 'E19' '  +0.60%  '
Set-TextValue 'D20' '15.18'
Set-TextValue 'E20' '  +1.92%  '
Set-TextValue 'D21' '26.951.10'
Set-TextValue 'E21' '  +3.37%  '
Set-TextValue 'D22' '5.197'
Set-TextValue 'E22' '  +4.43%  '
Set-TextValue 'D23' '10.97'
Set-TextValue 'E23' '  -0.89%  '
Set-TextValue 'D24' '2.279.93'
Set-TextValue 'E24' '  +16.46%  '
Set-TextValue 'E25' '  +1.61%  '
Set-TextValue 'D26' '1.886'
Set-TextValue 'E26' '  -0.85%  '
Set-TextValue 'D27' '18.33'
Set-TextValue 'E27' '  +1.42%  '
Set-TextValue 'D28' '5.242'
Set-TextValue 'E28' '  +4.51%  '
Set-TextValue 'D29' '1.928'
Set-TextValue 'E29' '  +15.81%  '
Set-TextValue 'D30' '114.77'
Set-TextValue 'E30' '  +0.59%  '
Set-TextValue 'D31' '0.08956'
Set-TextValue 'E31' '  +1.08%  '
Set-TextValue 'D32' '0.7429'
Set-TextValue 'E32' '  +3.16%  '
Set-TextValue 'D33' '1.158'
Set-TextValue 'E33' '  +5.12%  '
Set-TextValue 'D34' '4.434'
Set-TextValue 'E34' '  +3.21%  '
Set-TextValue 'D35' '2.806'
Set-TextValue 'E35' '  +4.44%  '
Set-TextValue 'E36' '  +0.65%  '
Set-TextValue 'D37' '1.124'
Set-TextValue 'E37' '  +5.62%  '
Set-TextValue 'D38' '0.05233'
Set-TextValue 'E38' '  +2.89%  '
Set-TextValue 'D39' '0.01923'
Set-TextValue 'E39' '  +2.68%  '
Set-TextValue 'D40' '0.5090'
Set-TextValue 'E40' '  +4.33%  '
Set-TextValue 'D41' '2.746'
Set-TextValue 'E41' '  +10.91%  '
Set-TextValue 'E42' '  +3.35%  '
Set-TextValue 'D43' '6.498'
Set-TextValue 'E43' '  +6.02%  '
Set-TextValue 'D44' '8.321'
Set-TextValue 'E44' '  +4.36%  '
Set-TextValue 'D45' '107.54'
Set-TextValue 'E45' '  +3.20%  '
Set-TextValue 'D46' '10.39'
Set-TextValue 'E46' '  +3.76%  '
Set-TextValue 'E47' '  +0.64%  '
Set-TextValue 'B48' 'Decentraland'
Set-TextValue 'C48' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D48' '0.4576'
Set-TextValue 'E48' '  +2.66%  '
Set-TextValue 'B49' 'NEARProtocol'
Set-TextValue 'C49' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D49' '1.650'
Set-TextValue 'E49' '  +5.50%  '
Set-TextValue 'B50' 'Cronos'
Set-TextValue 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D50' '0.06269'
Set-TextValue 'E50' '  +1.54%  '
Set-TextValue 'D51' '1.811'
Set-TextValue 'E51' '  +6.16%  '
